# wip: shadow --> boxShadow
#
# 1. On "Layout and Theme vars" (sheet1): mark the fully-completed rows
#    (shadow/boxShadow, textAlign, textAlignLast, textDecoration*, textTransform,
#    textWrap, userSelect, zIndex) with the "Good" (green) cell style, the same
#    style already used on finished rows like E58:E65 / row 66 / row 72 / row 74.
# 2. On "Changelog" (sheet2): split the single "Changes" list into two
#    sections - "Changes (layout properties)" and "Changes (theme variables)" -
#    and log the new "shadow --> boxShadow" theme-variable rename at the end.

$wb = $excel.ActiveWorkbook

$wsLayout = $wb.Worksheets.Item(1)
$wsChangelog = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# Sheet 1 - "Layout and Theme vars": flip newly-finished rows to the "Good"
# style. Use a format-only paste from an already-"Good" cell so the existing
# shared cellXf (style index 4) is reused instead of a new one being created.
# ---------------------------------------------------------------------------
$xlPasteFormats = -4122

$wsLayout.Range("E58").Copy()
$wsLayout.Range("A55:D65").PasteSpecial($xlPasteFormats)

$wsLayout.Range("E58").Copy()
$wsLayout.Range("E55:E57").PasteSpecial($xlPasteFormats)

$wsLayout.Range("E58").Copy()
$wsLayout.Range("A67:D67").PasteSpecial($xlPasteFormats)

$wsLayout.Range("E58").Copy()
$wsLayout.Range("A71:D71").PasteSpecial($xlPasteFormats)

$wsLayout.Range("E58").Copy()
$wsLayout.Range("A73:D73").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------------
# Sheet 2 - "Changelog": rename the header, insert a second section header,
# and append the new changelog line.
# ---------------------------------------------------------------------------

# Existing "Changes" header now specifically covers the layout-property renames.
$wsChangelog.Range("A1").Value = "Changes (layout properties)"

# Make room for a second section, right above the theme-variable entries
# (old row 11, "offset-decoration --> textUnderlineOffset").
$wsChangelog.Rows("11:11").Insert()

# New section header, styled the same as the first one (bold "Changes" style).
$wsChangelog.Range("A11").Value = "Changes (theme variables)"
$wsChangelog.Range("A1").Copy()
$wsChangelog.Range("A11").PasteSpecial($xlPasteFormats)

# Log the shadow --> boxShadow theme-variable rename as the newest entry.
$wsChangelog.Range("A21").Value = "shadow --> boxShadow"

# ---------------------------------------------------------------------------
# Restore selections to match the author's final cursor positions, with
# "Layout and Theme vars" left as the active sheet/tab.
# ---------------------------------------------------------------------------
$wsChangelog.Range("A12").Select()
$wsLayout.Range("C68").Select()

Write-Output "Applied: shadow --> boxShadow changelog + completed-row styling"
